# Add data for 2022-12-14
#
# Applies the per-cell deltas for one additional day of violent-crime
# incidents (2022-12-14) to the 2022 totals column (I) -- and a couple
# of 2015 column (B) corrections -- across the Citywide Totals sheet,
# the By Neighborhood summary sheet, and the affected per-neighborhood
# sheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 9).Value = 6988
$ws.Cells.Item(3, 9).Value = 7235
$ws.Cells.Item(4, 2).Value = 1667
$ws.Cells.Item(4, 9).Value = 1662
$ws.Cells.Item(5, 9).Value = 681
$ws.Cells.Item(6, 9).Value = 8528
$ws.Cells.Item(7, 2).Value = 23299
$ws.Cells.Item(7, 9).Value = 25094

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6, 9).Value = 183
$ws.Cells.Item(7, 9).Value = 787
$ws.Cells.Item(8, 9).Value = 1493
$ws.Cells.Item(9, 9).Value = 133
$ws.Cells.Item(11, 9).Value = 379
$ws.Cells.Item(18, 9).Value = 196
$ws.Cells.Item(19, 9).Value = 703
$ws.Cells.Item(29, 9).Value = 1493
$ws.Cells.Item(31, 9).Value = 254
$ws.Cells.Item(33, 9).Value = 1105
$ws.Cells.Item(34, 9).Value = 113
$ws.Cells.Item(36, 9).Value = 342
$ws.Cells.Item(37, 9).Value = 775
$ws.Cells.Item(43, 9).Value = 214
$ws.Cells.Item(44, 9).Value = 191
$ws.Cells.Item(47, 9).Value = 182
$ws.Cells.Item(48, 9).Value = 319
$ws.Cells.Item(49, 9).Value = 166
$ws.Cells.Item(52, 9).Value = 569
$ws.Cells.Item(54, 9).Value = 491
$ws.Cells.Item(57, 9).Value = 103
$ws.Cells.Item(63, 2).Value = 371
$ws.Cells.Item(63, 9).Value = 78
$ws.Cells.Item(64, 9).Value = 199
$ws.Cells.Item(65, 9).Value = 584
$ws.Cells.Item(67, 9).Value = 950
$ws.Cells.Item(69, 9).Value = 54
$ws.Cells.Item(77, 9).Value = 153
$ws.Cells.Item(78, 9).Value = 334
$ws.Cells.Item(79, 9).Value = 720
$ws.Cells.Item(83, 9).Value = 541
$ws.Cells.Item(85, 9).Value = 1118
$ws.Cells.Item(88, 9).Value = 232
$ws.Cells.Item(89, 9).Value = 295
$ws.Cells.Item(91, 9).Value = 263
$ws.Cells.Item(98, 9).Value = 183
$ws.Cells.Item(101, 2).Value = 23299
$ws.Cells.Item(101, 9).Value = 25094

# Sheet 3: South Shore
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 9).Value = 421
$ws.Cells.Item(6, 9).Value = 293
$ws.Cells.Item(7, 9).Value = 1118

# Sheet 4: Norwood Park
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 9).Value = 21
$ws.Cells.Item(7, 9).Value = 54

# Sheet 5: Little Village
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 9).Value = 141
$ws.Cells.Item(7, 9).Value = 569

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(6, 9).Value = 102
$ws.Cells.Item(7, 9).Value = 379

# Sheet 7: Austin
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 9).Value = 440
$ws.Cells.Item(3, 9).Value = 432
$ws.Cells.Item(6, 9).Value = 481
$ws.Cells.Item(7, 9).Value = 1493

# Sheet 9: Auburn Gresham
$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2, 9).Value = 259
$ws.Cells.Item(3, 9).Value = 239
$ws.Cells.Item(7, 9).Value = 787

# Sheet 10: Uptown
$ws = $wb.Worksheets.Item(10)
$ws.Cells.Item(5, 9).Value = 9
$ws.Cells.Item(7, 9).Value = 295

# Sheet 14: Grand Crossing
$ws = $wb.Worksheets.Item(14)
$ws.Cells.Item(2, 9).Value = 230
$ws.Cells.Item(4, 9).Value = 38
$ws.Cells.Item(7, 9).Value = 775

# Sheet 16: North Lawndale
$ws = $wb.Worksheets.Item(16)
$ws.Cells.Item(3, 9).Value = 355
$ws.Cells.Item(6, 9).Value = 283
$ws.Cells.Item(7, 9).Value = 950

# Sheet 17: Gage Park
$ws = $wb.Worksheets.Item(17)
$ws.Cells.Item(3, 9).Value = 65
$ws.Cells.Item(7, 9).Value = 254

# Sheet 19: New City
$ws = $wb.Worksheets.Item(19)
$ws.Cells.Item(4, 9).Value = 24
$ws.Cells.Item(7, 9).Value = 584

# Sheet 20: South Chicago
$ws = $wb.Worksheets.Item(20)
$ws.Cells.Item(2, 9).Value = 180
$ws.Cells.Item(5, 9).Value = 23
$ws.Cells.Item(7, 9).Value = 541

# Sheet 22: Garfield Park
$ws = $wb.Worksheets.Item(22)
$ws.Cells.Item(2, 9).Value = 249
$ws.Cells.Item(3, 9).Value = 407
$ws.Cells.Item(6, 9).Value = 355
$ws.Cells.Item(7, 9).Value = 1105

# Sheet 23: Lincoln Park
$ws = $wb.Worksheets.Item(23)
$ws.Cells.Item(3, 9).Value = 18
$ws.Cells.Item(6, 9).Value = 97
$ws.Cells.Item(7, 9).Value = 166

# Sheet 24: Loop
$ws = $wb.Worksheets.Item(24)
$ws.Cells.Item(6, 9).Value = 237
$ws.Cells.Item(7, 9).Value = 491

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item(25)
$ws.Cells.Item(2, 9).Value = 443
$ws.Cells.Item(4, 9).Value = 80
$ws.Cells.Item(6, 9).Value = 411
$ws.Cells.Item(7, 9).Value = 1493

# Sheet 26: Chatham
$ws = $wb.Worksheets.Item(26)
$ws.Cells.Item(2, 9).Value = 227
$ws.Cells.Item(7, 9).Value = 703

# Sheet 27: Irving Park
$ws = $wb.Worksheets.Item(27)
$ws.Cells.Item(2, 9).Value = 64
$ws.Cells.Item(3, 9).Value = 54
$ws.Cells.Item(7, 9).Value = 191

# Sheet 28: Lake View
$ws = $wb.Worksheets.Item(28)
$ws.Cells.Item(3, 9).Value = 61
$ws.Cells.Item(7, 9).Value = 319

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item(30)
$ws.Cells.Item(3, 9).Value = 44
$ws.Cells.Item(7, 9).Value = 183

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item(35)
$ws.Cells.Item(6, 9).Value = 122
$ws.Cells.Item(7, 9).Value = 334

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item(40)
$ws.Cells.Item(6, 9).Value = 70
$ws.Cells.Item(7, 9).Value = 263

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item(42)
$ws.Cells.Item(2, 9).Value = 210
$ws.Cells.Item(3, 9).Value = 237
$ws.Cells.Item(4, 9).Value = 42
$ws.Cells.Item(7, 9).Value = 720

# Sheet 43: Near South Side
$ws = $wb.Worksheets.Item(43)
$ws.Cells.Item(3, 9).Value = 57
$ws.Cells.Item(7, 9).Value = 199

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item(45)
$ws.Cells.Item(3, 9).Value = 44
$ws.Cells.Item(7, 9).Value = 196

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item(47)
$ws.Cells.Item(2, 9).Value = 97
$ws.Cells.Item(7, 9).Value = 342

# Sheet 50: Garfield Ridge
$ws = $wb.Worksheets.Item(50)
$ws.Cells.Item(2, 9).Value = 46
$ws.Cells.Item(7, 9).Value = 113

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item(53)
$ws.Cells.Item(2, 9).Value = 43
$ws.Cells.Item(7, 9).Value = 182

# Sheet 55: Wicker Park
$ws = $wb.Worksheets.Item(55)
$ws.Cells.Item(6, 9).Value = 119
$ws.Cells.Item(7, 9).Value = 183

# Sheet 61: Avalon Park
$ws = $wb.Worksheets.Item(61)
$ws.Cells.Item(3, 9).Value = 44
$ws.Cells.Item(7, 9).Value = 133

# Sheet 68: United Center
$ws = $wb.Worksheets.Item(68)
$ws.Cells.Item(6, 9).Value = 75
$ws.Cells.Item(7, 9).Value = 232

# Sheet 77: Mckinley Park
$ws = $wb.Worksheets.Item(77)
$ws.Cells.Item(2, 9).Value = 39
$ws.Cells.Item(7, 9).Value = 103

# Sheet 79: Hyde Park
$ws = $wb.Worksheets.Item(79)
$ws.Cells.Item(2, 9).Value = 44
$ws.Cells.Item(7, 9).Value = 214

# Sheet 84: Riverdale
$ws = $wb.Worksheets.Item(84)
$ws.Cells.Item(2, 9).Value = 51
$ws.Cells.Item(4, 9).Value = 7
$ws.Cells.Item(7, 9).Value = 153

